$wb = $excel.ActiveWorkbook

# --- Worksheets -------------------------------------------------------
$ws1 = $wb.Worksheets.Item("InstrumentAttribute")

# --- Update revenue figures on InstrumentAttribute (rows 23-26) -------
$ws1.Range("H23").Value = -1000
$ws1.Range("J23").Value = -1000

$ws1.Range("H24").Value = 500
$ws1.Range("J24").Value = 6000

$ws1.Range("H25").Value = 1500
$ws1.Range("J25").Value = 1500

$ws1.Range("H26").Value = 250
$ws1.Range("J26").Value = 3000

# --- Switch the active / selected sheet from xxProduct to ----------
# --- InstrumentAttribute, resetting its selection to the default A1 -
$ws1.Activate() | Out-Null
$ws1.Range("A1").Select() | Out-Null
